$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C (plain text, non-numeric) can be set directly.
# Columns D/E contain numeric-looking / percentage-looking text that must
# stay stored as literal text (matching the original inlineStr cells), so we
# force a text number-format before assignment, then clear the format again
# (ClearFormats) so no stray style index is left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '326.61'
$ws.Range('D2').ClearFormats()

$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '-1.04%'
$ws.Range('E2').ClearFormats()

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '44.24'
$ws.Range('D3').ClearFormats()

$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '0.92%'
$ws.Range('E3').ClearFormats()

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '5.545'
$ws.Range('D4').ClearFormats()

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '-4.02%'
$ws.Range('E4').ClearFormats()

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.08079'
$ws.Range('D5').ClearFormats()

$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-2.80%'
$ws.Range('E5').ClearFormats()

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '8.703'
$ws.Range('D6').ClearFormats()

$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '-0.90%'
$ws.Range('E6').ClearFormats()

$ws.Range('B7').Value = 'FTXToken'

$ws.Range('C7').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.915'
$ws.Range('D7').ClearFormats()

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-2.29%'
$ws.Range('E7').ClearFormats()

$ws.Range('B8').Value = 'GateToken'

$ws.Range('C8').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '4.332'
$ws.Range('D8').ClearFormats()

$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-3.71%'
$ws.Range('E8').ClearFormats()

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.704'
$ws.Range('D9').ClearFormats()

$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-7.45%'
$ws.Range('E9').ClearFormats()

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9489'
$ws.Range('D10').ClearFormats()

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '1.41%'
$ws.Range('E10').ClearFormats()

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.1181'
$ws.Range('D11').ClearFormats()

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-5.47%'
$ws.Range('E11').ClearFormats()

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.1895'
$ws.Range('D12').ClearFormats()

$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '-2.32%'
$ws.Range('E12').ClearFormats()

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.1008'
$ws.Range('D13').ClearFormats()

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '6.10%'
$ws.Range('E13').ClearFormats()

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.04183'
$ws.Range('D14').ClearFormats()

$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '5.16%'
$ws.Range('E14').ClearFormats()

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.1066'
$ws.Range('D15').ClearFormats()

$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '0.18%'
$ws.Range('E15').ClearFormats()

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.001279'
$ws.Range('D16').ClearFormats()

$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-1.67%'
$ws.Range('E16').ClearFormats()

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.006031'
$ws.Range('D17').ClearFormats()

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '1.48%'
$ws.Range('E17').ClearFormats()

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.600'
$ws.Range('D18').ClearFormats()

$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '2.13%'
$ws.Range('E18').ClearFormats()

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.64%'
$ws.Range('E19').ClearFormats()

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.409'
$ws.Range('D20').ClearFormats()

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '-7.19%'
$ws.Range('E20').ClearFormats()

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.1373'
$ws.Range('D21').ClearFormats()

$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '0.18%'
$ws.Range('E21').ClearFormats()

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '3.63%'
$ws.Range('E22').ClearFormats()

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04254'
$ws.Range('D23').ClearFormats()

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '-3.30%'
$ws.Range('E23').ClearFormats()

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001236'
$ws.Range('D24').ClearFormats()

$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-1.49%'
$ws.Range('E24').ClearFormats()

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004601'
$ws.Range('D25').ClearFormats()

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '6.06%'
$ws.Range('E25').ClearFormats()

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001234'
$ws.Range('D26').ClearFormats()

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0003999'
$ws.Range('D27').ClearFormats()

$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '0.15%'
$ws.Range('E27').ClearFormats()

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02660'
$ws.Range('D39').ClearFormats()

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '-4.96%'
$ws.Range('E39').ClearFormats()

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.05550'
$ws.Range('D40').ClearFormats()

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '-0.76%'
$ws.Range('E40').ClearFormats()

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '24.96%'
$ws.Range('E41').ClearFormats()

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.007690'
$ws.Range('D42').ClearFormats()

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '-2.98%'
$ws.Range('E42').ClearFormats()

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1394'
$ws.Range('D43').ClearFormats()

$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-2.09%'
$ws.Range('E43').ClearFormats()

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.002061'
$ws.Range('D44').ClearFormats()

$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-4.42%'
$ws.Range('E44').ClearFormats()

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.008695'
$ws.Range('D45').ClearFormats()

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-12.39%'
$ws.Range('E45').ClearFormats()

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00007131'
$ws.Range('D46').ClearFormats()

$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '-1.12%'
$ws.Range('E46').ClearFormats()

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.00000000752'
$ws.Range('D47').ClearFormats()

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '0.22%'
$ws.Range('E47').ClearFormats()

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.003435'
$ws.Range('D48').ClearFormats()

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '-14.80%'
$ws.Range('E48').ClearFormats()

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.002276'
$ws.Range('D49').ClearFormats()

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '-0.19%'
$ws.Range('E49').ClearFormats()

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.00002107'
$ws.Range('D50').ClearFormats()

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.22%'
$ws.Range('E50').ClearFormats()

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0002006'
$ws.Range('D51').ClearFormats()

$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '0.22%'
$ws.Range('E51').ClearFormats()
